$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.432.22'
$ws.Range('E2').Value = '  -0.95%  '
$ws.Range('D3').Value = '1.884.44'
$ws.Range('E3').Value = '  -0.97%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.45'
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4826'
$ws.Range('E7').Value = '  -2.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2885'
$ws.Range('E8').Value = '  -2.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06594'
$ws.Range('E9').Value = '  -1.88%  '
$ws.Range('D10').Value = '1.892.42'
$ws.Range('E10').Value = '  -0.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '16.89'
$ws.Range('E11').Value = '  -0.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07394'
$ws.Range('E12').Value = '  +0.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.178'
$ws.Range('E13').Value = '  +0.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.57'
$ws.Range('E14').Value = '  +0.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6600'
$ws.Range('E15').Value = '  -1.26%  '
$ws.Range('D16').Value = '30.413.31'
$ws.Range('E16').Value = '  -0.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.53'
$ws.Range('E17').Value = '  +0.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007751'
$ws.Range('E18').Value = '  -1.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').Value = '  -0.10%  '
$ws.Range('D20').Value = '2.149.91'
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.376'
$ws.Range('E21').Value = '  +0.70%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.006'
$ws.Range('E22').Value = '  +0.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '219.71'
$ws.Range('E23').Value = '  +15.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.171'
$ws.Range('E24').Value = '  -1.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.363'
$ws.Range('E25').Value = '  -2.33%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.33'
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.63'
$ws.Range('E27').Value = '  +0.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.935'
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.439'
$ws.Range('E29').Value = '  -2.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.328'
$ws.Range('E30').Value = '  -2.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09182'
$ws.Range('E31').Value = '  +0.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.042'
$ws.Range('E32').Value = '  +0.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05069'
$ws.Range('E33').Value = '  -3.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7527'
$ws.Range('E34').Value = '  +1.55%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.153'
$ws.Range('E35').Value = '  +4.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.713'
$ws.Range('E36').Value = '  -0.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01877'
$ws.Range('E37').Value = '  +2.79%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.644'
$ws.Range('E38').Value = '  -2.46%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9191'
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.087'
$ws.Range('E40').Value = '  +0.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.000'
$ws.Range('E41').Value = '  +1.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '106.71'
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4332'
$ws.Range('E43').Value = '  -2.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.002'
$ws.Range('E44').Value = '  +0.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.613'
$ws.Range('E45').Value = '  +1.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1332'
$ws.Range('E46').Value = '  -3.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.583'
$ws.Range('E47').Value = '  +10.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '65.05'
$ws.Range('E48').Value = '  -12.75%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.941'
$ws.Range('E49').Value = '  -1.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.36'
$ws.Range('E50').Value = '  -2.97%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05720'
$ws.Range('E51').Value = '  -2.24%  '
